# Populate the "RunMode" column (C) on the LoginTestdata sheet with
# pass/skip flags for each test-data row, and leave the selection on B2
# (matching the author's final cursor position when they saved the file).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginTestdata")

$ws.Range("C2").Value = "y"
$ws.Range("C3").Value = "N"
$ws.Range("C4").Value = "y"
$ws.Range("C5").Value = "y"

$ws.Activate()
$ws.Range("B2").Select()
